# Price update for 2026-02-07
# Appends one new tracking row (row 38) to Sheet 1 with the latest scrape:
#   Date = 2026-02-07, Price = 5250000, Discount = 70, Incredible = 0
#
# All existing rows store every value as text (shared strings), even the
# purely-numeric ones (e.g. "50", "0", "8749000"). Assigning a plain string
# to .Value would let Excel auto-coerce numeric-looking / date-looking text
# into a real number/date, so we force the new cells to Text format first,
# write the values, then drop the explicit format again so the cells end up
# indistinguishable from the rest of the (unstyled) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$target = $ws.Range("A" + $row + ":D" + $row)
$target.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2026-02-07"
$ws.Range("B" + $row).Value = "5250000"
$ws.Range("C" + $row).Value = "70"
$ws.Range("D" + $row).Value = "0"

# Restore the default style so the new cells don't carry a leftover
# explicit "Text" number format that the rest of the sheet doesn't have.
$target.Style = "Normal"
